# Updated queries for C3DC first half testcases.
#
# The six/seven DuckDB-style SQL statements stored on Sheet1 (columns B & C,
# rows 2-7) joined tables using the non-namespaced "id" column
# (std.id / prt.id) together with oddly-named foreign-key columns
# ("study.id" / "participant.id"). The source tables were renamed so the
# natural keys are now "study_id" / "participant_id" everywhere - update
# every occurrence of the old join predicates to the new column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Fix-Sql([string]$text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

# Walk every used cell and rewrite any SQL text that still references the
# old join columns - this hits the StatQuery (C2) and each TabQuery (B2:B7)
# cell without hard-coding row numbers.
$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value()
    if (($val -ne $null) -and ($val -is [string]) -and $val.Contains('df_study')) {
        $newVal = Fix-Sql $val
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# The author's last interactive selection when the workbook was saved moved
# from C5 to C7 (scrolled down one row further to the Survival tab's query).
$ws.Range("C7").Select()
